$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell A26: "conda install --channel https://conda.anac.org/menpo opencv3\n"
#        -> "conda install --channel https://conda.anac.org/menpo opencv3" (trailing blank line removed)
$ws.Range("A26").Value = "conda install --channel https://conda.anac.org/menpo opencv3"

# Cell A27: "`nimport sys`n##sys.path.append('/usr/local/lib/python2.7/site-packages')`nsys.path.append('/home/kshiba/conda/lib/python3.6/site-packages')`nimport cv2`n"
#        -> "import sys`n##sys.path.append('/usr/local/lib/python2.7/site-packages')`nsys.path.append('/home/kshiba/conda/lib/python3.6/site-packages')`nimport cv2" (leading & trailing blank lines removed)
$ws.Range("A27").Value = "import sys`n##sys.path.append('/usr/local/lib/python2.7/site-packages')`nsys.path.append('/home/kshiba/conda/lib/python3.6/site-packages')`nimport cv2"

# Row 27 is now two lines shorter, so its (wrap-text driven) height shrinks from 64.15 to 43.25
$ws.Rows("27").RowHeight = 43.25

# Move the selection from C27 to A29, and scroll the view down
$ws.Range("A29").Select()
